# Applies the contract amendment described by the commit diff:
#  - swap the MUTUÁRIO(A) party from "C. V. LEITE DE LIMA" to "AMBEV S.A."
#    (CNPJ, phone, e-mail, address and contract number all change too)
#  - switch repayment from "parcelado em 5x" to a single installment
#  - update the delivery / first-installment date
#  - bump the interest rate 25.0% -> 25.8%
#  - remove installments 2 through 5 from the payment schedule
#  - update the contract date and the MUTUÁRIO signature block

$d = $word.ActiveDocument

# 1) Opening paragraph: replace the whole MUTUÁRIO(A) description in one shot
#    (this is the only occurrence of this long run, so Find/Replace is safe).
$old1 = "C. V. LEITE DE LIMA,inscrito no CNPJ: 11004949000186, telefone nº (95) 3224-9265, e-mail: ce@gmail.com, residente e domiciliado na  com sede social a Rua:  AVENIDA PRESIDENTE JUSCELINO KUBITSCHECK, nº 1032 , NOSSA SENHORA APARECIDA, CEP: 69306295 e SÃO PAULO/RR, doravante denominado(a) MUTUÁRIO(A), o presente mútuo, contrato nº 5 mediante as seguintes cláusulas:"
$new1 = "AMBEV S.A.,inscrito no CNPJ: 07526557011659, telefone nº (19) 3313-5680, e-mail: opobrigaces@ambev.com.br, residente e domiciliado na  com sede social a Rua:  AV CONSTANTINO NERY, nº 2575 ANDAR 01 AO 8 E 9 ANDAR CONJ 902 E 16 ANDAR, FLORES, CEP: 69058795 e São Paulo/AM, doravante denominado(a) MUTUÁRIO(A), o presente mútuo, contrato nº 2 mediante as seguintes cláusulas:"
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2) 1.2. payment method: parcelado em 5x -> uma única parcela
$d.Content.Find.Execute(
    "1.2. A quantia será repassada ao(à) MUTUÁRIO(A) mediante, por meio parcelado em 5x.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "1.2. A quantia será repassada ao(à) MUTUÁRIO(A) mediante, por meio uma única parcela.", 2)

# 3) 1.3. delivery date: 03/12/2025 -> 13/02/2026
$d.Content.Find.Execute(
    "1.3. O(A) MUTUANTE entregará a quantia ao(à) MUTUÁRIO(A) no ato de assinatura deste instrumento OU em 03/12/2025. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "1.3. O(A) MUTUANTE entregará a quantia ao(à) MUTUÁRIO(A) no ato de assinatura deste instrumento OU em 13/02/2026. ", 2)

# 4) Payment schedule heading: Parcelado em 5x -> Parcela única
$d.Content.Find.Execute("Parcelado em 5x", $false, $false, $false, $false, $false, $true, 1, $false, "Parcela única", 2)

# 5) First installment date updates to match the new single-payment date
#    (replace only the date text itself so the bullet/tab runs are untouched)
$d.Content.Find.Execute("03/12/2025 ", $false, $false, $false, $false, $false, $true, 1, $false, "13/02/2026 ", 2)

# 6) Remove installments 2ª-5ª (whole paragraphs, including the paragraph mark)
$bulletsToRemove = @("2ª 03/01/2026", "3ª 03/02/2026", "4ª 03/03/2026", "5ª 03/04/2026")
foreach ($bullet in $bulletsToRemove) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ("*" + $bullet + "*")) {
            $p.Range.Delete()
            break
        }
    }
}

# 7) Interest rate: 25.0% -> 25.8%
$d.Content.Find.Execute("25.0%", $false, $false, $false, $false, $false, $true, 1, $false, "25.8%", 2)

# 8) Contract date at the bottom: São Paulo, 02/11/2025. -> São Paulo, 11/01/2026.
$d.Content.Find.Execute("São Paulo, 02/11/2025.", $false, $false, $false, $false, $false, $true, 1, $false, "São Paulo, 11/01/2026.", 2)

# 9) MUTUÁRIO signature block: name + CNPJ
$d.Content.Find.Execute("C. V. LEITE DE LIMA", $false, $false, $false, $false, $false, $true, 1, $false, "AMBEV S.A.", 2)
$d.Content.Find.Execute("11004949000186", $false, $false, $false, $false, $false, $true, 1, $false, "07526557011659", 2)

Write-Output "done"
